$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "2025/12/03 01:48"
$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"
$ws.Range("F8").Value = "-"
$ws.Range("G8").Value = "-"
